# Weekly refresh of "Fruto del paraíso" hortaliza data:
# the rows 2-8 are re-ordered (and some fields updated) while
# columns A,B,C,E,F,G,H,O,R stay identical across all rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values, keyed by destination row, for columns D,I,J,K,L,M,N,P,Q
$rows = @{
    2 = @{ D = 44293; I = "Primera";  J = 10; K = 25000; L = 25000; M = 25000; N = "`$/caja 15 kilos empedrada"; P = 1667; Q = 15 }
    3 = @{ D = 44285; I = "Primera";  J = 20; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
    4 = @{ D = 44315; I = "Especial"; J = 10; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 }
    5 = @{ D = 44315; I = "Primera";  J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos granel";    P = 1000; Q = 15 }
    6 = @{ D = 44313; I = "Primera";  J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos empedrada"; P = 1000; Q = 15 }
    7 = @{ D = 44313; I = "Primera";  J = 20; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 }
    8 = @{ D = 44280; I = "Primera";  J = 30; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value2  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 9).Value   = $vals.I   # I: Calidad
    $ws.Cells.Item($r, 10).Value2 = $vals.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value2 = $vals.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $vals.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value  = $vals.N   # N: Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value2 = $vals.P   # P: Precio $/Kg
    $ws.Cells.Item($r, 17).Value2 = $vals.Q   # Q: Kg o Unidades
}
